$wb = $excel.ActiveWorkbook

# --- Sheet 1: _validation_data ---
$ws1 = $wb.Worksheets.Item("_validation_data")

# Move column C (time units) into column B (overwriting vendor names) for rows 1-5
$ws1.Range("B1").Value = $ws1.Range("C1").Value2
$ws1.Range("B2").Value = $ws1.Range("C2").Value2
$ws1.Range("B3").Value = $ws1.Range("C3").Value2
$ws1.Range("B4").Value = $ws1.Range("C4").Value2
$ws1.Range("B5").Value = $ws1.Range("C5").Value2

# Clear leftover vendor names in column B for rows 6-9 (no longer used)
$ws1.Range("B6:B9").ClearContents()

# Delete now-redundant column C
$ws1.Columns.Item(3).Delete()

# Insert "Custom" as new row 6
$ws1.Rows.Item(6).Insert()
$ws1.Range("A6").Value = "Custom"

Write-Host "Sheet1 done"

# --- Sheet 2: Non-Standard Value ---
$ws2 = $wb.Worksheets.Item("Non-Standard Value")

# Remove the "preparation_instrument_vendor" rows (bottom-up so row numbers stay valid)
$ws2.Rows.Item(11).Delete()
$ws2.Rows.Item(9).Delete()
$ws2.Rows.Item(7).Delete()
$ws2.Rows.Item(5).Delete()
$ws2.Rows.Item(3).Delete()

# Re-apply the AutoFilter over the new, smaller range
$ws2.AutoFilterMode = $false
$ws2.Range("A1:E6").AutoFilter() | Out-Null

# Fix up the data validation list reference (sheet1 grew by one row)
$dv2 = $ws2.Range("D2:D6").Validation
$dv2.Modify($dv2.Type, $dv2.AlertStyle, $dv2.Operator, "_validation_data!`$A`$1:`$A`$19")

Write-Host "Sheet2 done"

# --- Sheet 3: Missing Required Value ---
$ws3 = $wb.Worksheets.Item("Missing Required Value")
$dv3 = $ws3.Range("C3:C15").Validation
$dv3.Modify($dv3.Type, $dv3.AlertStyle, $dv3.Operator, "_validation_data!`$B`$1:`$B`$5")

Write-Host "Sheet3 done"

# --- Workbook level: fix up the _FilterDatabase defined name for sheet 2 ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $nm = $wb.Names.Item($i)
  if ($nm.Name -eq "Non-Standard Value!_FilterDatabase") {
    $nm.RefersTo = "='Non-Standard Value'!`$A`$1:`$E`$6"
  }
}

Write-Host "Workbook names done"
